$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Using GPU" checkbox -> TRUE
$ws.Range("B2").Value = $true

# "Bias(V)" -> 1
$ws.Range("B9").Value = 1

# "Plot band structure" checkbox -> TRUE
$ws.Range("B10").Value = $true

# Data table row 14 ("o"): Length(# of unit cell) 808 -> 1000, Gap Open (eV) 0.05 -> 0.1
$ws.Range("F14").Value = 1000
$ws.Range("I14").Value = 0.1

# Data table row 15 ("x"): Gap Open (eV) 0.05 -> 0.1
$ws.Range("I15").Value = 0.1

# Data table row 16 ("x"): Width(# of sub unit cell) 1 -> 0, Gap Open (eV) 0.05 -> 0.1
$ws.Range("E16").Value = 0
$ws.Range("I16").Value = 0.1

# Fill Width(nm)/Length(nm) formulas down J15:J16 / K15:K16 as shared formulas
# (matches how Excel stores a formula typed once and filled down a column)
$ws.Range("J15:J16").Formula = "=E15*0.246*3^0.5/2"
$ws.Range("K15:K16").Formula = "=F15*0.246*3"

# Move the active selection to B3 (matches the saved UI state in the target file)
$ws.Range("B3").Select()
